$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the glucose-range bucket labels from "XX - YY" to "XX até YY".
# (B1 "Pessoas" text itself is unchanged; its shared-string index shifts
# automatically as a natural side effect of the table being rewritten.)
$ws.Range("A2").Value = "56 até 74"
$ws.Range("A3").Value = "74 até 92"
$ws.Range("A4").Value = "93 até 110"
$ws.Range("A5").Value = "111 até 128"
$ws.Range("A6").Value = "129 até 146"
$ws.Range("A7").Value = "147 até 164"
$ws.Range("A8").Value = "165 até 182"
$ws.Range("A9").Value = "183 até 199"

# Move the active selection to match the saved view state in the target file.
$ws.Range("J11").Select() | Out-Null
